$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 37. This shifts every row from 37 downward by one,
# which reproduces the row-37..56 -> 38..57 shift seen in the diff in one shot.
$ws.Rows.Item(37).Insert()

# Row 38 (the old row 37, now shifted down): extend the existing TODO text first
# so it keeps claiming the earlier shared-string slot, matching the author's edit order.
$ws.Cells.Item(38, 4).Value2 = "Todo - complete javadocs for user; figure out where to test cascading deletes; complete javadocs for story; refactor tests to use .equals"
$ws.Rows.Item(38).RowHeight = 30

# New row 37: brand new TODO note (style matches the wrap-text "Task" column style
# already used by row 38/neighbours, carried over automatically by Insert()).
$ws.Cells.Item(37, 4).Value2 = "TODO - Add unit tests to thoroughly test each method in your DAOs. Be sure that you are testing insert and delete of the associated records."
$ws.Rows.Item(37).RowHeight = 30

# Update the view: scroll position + active cell/selection.
$ws.Range("G37").Select()
